$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "5721736"
$ws.Cells.Item(3, 2).Value = "Super Silk Taschentücher Würfelbox"
$ws.Cells.Item(3, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/super-silk-taschentuecher-wuerfelbox/p/5721736"
$ws.Cells.Item(3, 4).Value = "60ST"
$ws.Cells.Item(3, 5).Value = 10
$ws.Cells.Item(3, 6).Value = 4
$ws.Cells.Item(3, 7).Value = "Super Silk"
$ws.Cells.Item(3, 8).NumberFormat = "@"
$ws.Cells.Item(3, 8).Value = "2.10"
$ws.Cells.Item(3, 9).Value = "0.04/1ST"
$ws.Cells.Item(3, 11).NumberFormat = "@"
$ws.Cells.Item(3, 11).Value = "0.04"
$ws.Cells.Item(3, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Cells.Item(3, 14).Value = "Super Silk Taschentücher Würfelbox 2.10 Schweizer Franken"
$ws.Cells.Item(3, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "6568452"
$ws.Cells.Item(4, 2).Value = "Super Soft Premium Mandel feucht 4x  50ST"
$ws.Cells.Item(4, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/super-soft-premium-mandel-feucht/p/6568452"
$ws.Cells.Item(4, 4).Value = "4x 50ST"
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = 3.5
$ws.Cells.Item(4, 7).Value = "Super Soft"
$ws.Cells.Item(4, 8).NumberFormat = "@"
$ws.Cells.Item(4, 8).Value = "7.65"
$ws.Cells.Item(4, 9).Value = "0.04/1ST"
$ws.Cells.Item(4, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(4, 11).NumberFormat = "@"
$ws.Cells.Item(4, 11).Value = "0.04"
$ws.Cells.Item(4, 12).Value = "1ST"
$ws.Cells.Item(4, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Cells.Item(4, 14).Value = "Super Soft Premium Mandel feucht 4x  50ST 35% Aktion 7.65 Schweizer Franken statt 11.80 Schweizer Franken"
$ws.Cells.Item(4, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "6283679"
$ws.Cells.Item(5, 2).Value = "Oecoplan Toilettenpapier Camomille weiss 4-lagig 6 Rollen"
$ws.Cells.Item(5, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/oecoplan-toilettenpapier-camomille-weiss-4-lagig-6-rollen/p/6283679"
$ws.Cells.Item(5, 4).Value = "6Rol"
$ws.Cells.Item(5, 5).Value = 13
$ws.Cells.Item(5, 6).Value = 4
$ws.Cells.Item(5, 7).Value = "Coop"
$ws.Cells.Item(5, 8).NumberFormat = "@"
$ws.Cells.Item(5, 8).Value = "4.50"
$ws.Cells.Item(5, 9).Value = "0.75/1Rol"
$ws.Cells.Item(5, 10).Value = "Preis pro 1 Rolle"
$ws.Cells.Item(5, 11).NumberFormat = "@"
$ws.Cells.Item(5, 11).Value = "0.75"
$ws.Cells.Item(5, 12).Value = "1Rol"
$ws.Cells.Item(5, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Cells.Item(5, 14).Value = "Oecoplan Toilettenpapier Camomille weiss 4-lagig 6 Rollen 4.50 Schweizer Franken"
$ws.Cells.Item(5, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(6, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(7, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(8, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(9, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(10, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "6498157"
$ws.Cells.Item(11, 2).Value = "subito Haushaltspapier weiss 2 Rollen"
$ws.Cells.Item(11, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/subito-haushaltspapier-weiss-2-rollen/p/6498157"
$ws.Cells.Item(11, 4).Value = "100BLT"
$ws.Cells.Item(11, 5).Value = ""
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = "subito"
$ws.Cells.Item(11, 8).NumberFormat = "@"
$ws.Cells.Item(11, 8).Value = "2.75"
$ws.Cells.Item(11, 9).Value = ""
$ws.Cells.Item(11, 10).Value = ""
$ws.Cells.Item(11, 11).Value = ""
$ws.Cells.Item(11, 12).Value = ""
$ws.Cells.Item(11, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
$ws.Cells.Item(11, 14).Value = "subito Haushaltspapier weiss 2 Rollen 2.75 Schweizer Franken"
$ws.Cells.Item(11, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "6724076"
$ws.Cells.Item(12, 2).Value = "Oecoplan feuchtes Toilettenpapier Duckies natural 40 Stück"
$ws.Cells.Item(12, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/oecoplan-feuchtes-toilettenpapier-duckies-natural-40-stueck/p/6724076"
$ws.Cells.Item(12, 4).Value = "40ST"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 4.5
$ws.Cells.Item(12, 7).Value = "Duckies"
$ws.Cells.Item(12, 8).NumberFormat = "@"
$ws.Cells.Item(12, 8).Value = "2.95"
$ws.Cells.Item(12, 9).Value = "0.07/1ST"
$ws.Cells.Item(12, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(12, 11).NumberFormat = "@"
$ws.Cells.Item(12, 11).Value = "0.07"
$ws.Cells.Item(12, 12).Value = "1ST"
$ws.Cells.Item(12, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Cells.Item(12, 14).Value = "Oecoplan feuchtes Toilettenpapier Duckies natural 40 Stück 2.95 Schweizer Franken"
$ws.Cells.Item(12, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(13, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Value = "3180824"
$ws.Cells.Item(14, 2).Value = "Tempo Taschentücher Plus Aloe &amp; Kamille 12x9 Stück"
$ws.Cells.Item(14, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/tempo-taschentuecher-plus-aloe-kamille-12x9-stueck/p/3180824"
$ws.Cells.Item(14, 4).Value = "12ST"
$ws.Cells.Item(14, 5).Value = 7
$ws.Cells.Item(14, 6).Value = 4
$ws.Cells.Item(14, 8).NumberFormat = "@"
$ws.Cells.Item(14, 8).Value = "3.95"
$ws.Cells.Item(14, 9).Value = "0.33/1ST"
$ws.Cells.Item(14, 11).NumberFormat = "@"
$ws.Cells.Item(14, 11).Value = "0.33"
$ws.Cells.Item(14, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Cells.Item(14, 14).Value = "Tempo Taschentücher Plus Aloe &amp; Kamille 12x9 Stück 3.95 Schweizer Franken"
$ws.Cells.Item(14, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = "6711017"
$ws.Cells.Item(15, 2).Value = "Tempo feuchte Limited Edition"
$ws.Cells.Item(15, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/tempo-feuchte-limited-edition/p/6711017"
$ws.Cells.Item(15, 4).Value = "42ST"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 2.5
$ws.Cells.Item(15, 8).NumberFormat = "@"
$ws.Cells.Item(15, 8).Value = "3.45"
$ws.Cells.Item(15, 9).Value = "0.08/1ST"
$ws.Cells.Item(15, 11).NumberFormat = "@"
$ws.Cells.Item(15, 11).Value = "0.08"
$ws.Cells.Item(15, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Cells.Item(15, 14).Value = "Tempo feuchte Limited Edition 3.45 Schweizer Franken"
$ws.Cells.Item(15, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(16, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(17, 1).NumberFormat = "@"
$ws.Cells.Item(17, 1).Value = "6433417"
$ws.Cells.Item(17, 2).Value = "Plenty Fun Design extra Long"
$ws.Cells.Item(17, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/plenty-fun-design-extra-long/p/6433417"
$ws.Cells.Item(17, 4).Value = "144BLT"
$ws.Cells.Item(17, 5).Value = 6
$ws.Cells.Item(17, 6).Value = 4.5
$ws.Cells.Item(17, 7).Value = "Plenty"
$ws.Cells.Item(17, 8).NumberFormat = "@"
$ws.Cells.Item(17, 8).Value = "6.70"
$ws.Cells.Item(17, 14).Value = "Plenty Fun Design extra Long 6.70 Schweizer Franken"
$ws.Cells.Item(17, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(18, 1).NumberFormat = "@"
$ws.Cells.Item(18, 1).Value = "6638996"
$ws.Cells.Item(18, 2).Value = "Oecoplan Allzweck Papiertücher"
$ws.Cells.Item(18, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/oecoplan-allzweck-papiertuecher/p/6638996"
$ws.Cells.Item(18, 4).Value = "176BLT"
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 5
$ws.Cells.Item(18, 7).Value = "Coop"
$ws.Cells.Item(18, 8).NumberFormat = "@"
$ws.Cells.Item(18, 8).Value = "3.20"
$ws.Cells.Item(18, 14).Value = "Oecoplan Allzweck Papiertücher 3.20 Schweizer Franken"
$ws.Cells.Item(18, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(19, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(20, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(21, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(22, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(23, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(24, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(25, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(26, 1).NumberFormat = "@"
$ws.Cells.Item(26, 1).Value = "6283677"
$ws.Cells.Item(26, 2).Value = "Oecoplan Goldmelisse blau 3-lagig 32 Rollen"
$ws.Cells.Item(26, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/oecoplan-goldmelisse-blau-3-lagig-32-rollen/p/6283677"
$ws.Cells.Item(26, 4).Value = "32Rol"
$ws.Cells.Item(26, 6).Value = 5
$ws.Cells.Item(26, 7).Value = "Coop"
$ws.Cells.Item(26, 8).NumberFormat = "@"
$ws.Cells.Item(26, 8).Value = "14.80"
$ws.Cells.Item(26, 9).Value = "0.46/1Rol"
$ws.Cells.Item(26, 10).Value = "Preis pro 1 Rolle"
$ws.Cells.Item(26, 11).NumberFormat = "@"
$ws.Cells.Item(26, 11).Value = "0.46"
$ws.Cells.Item(26, 12).Value = "1Rol"
$ws.Cells.Item(26, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Cells.Item(26, 14).Value = "Oecoplan Goldmelisse blau 3-lagig 32 Rollen 30% Aktion 14.80 Schweizer Franken statt 21.20 Schweizer Franken"
$ws.Cells.Item(26, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(27, 1).NumberFormat = "@"
$ws.Cells.Item(27, 1).Value = "6996030"
$ws.Cells.Item(27, 2).Value = "Tela Viva Haushaltspapier 3-lagig 4 Rollen"
$ws.Cells.Item(27, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/tela-viva-haushaltspapier-3-lagig-4-rollen/p/6996030"
$ws.Cells.Item(27, 4).Value = "200BLT"
$ws.Cells.Item(27, 5).Value = 1
$ws.Cells.Item(27, 6).Value = 4
$ws.Cells.Item(27, 7).Value = "Tela"
$ws.Cells.Item(27, 8).NumberFormat = "@"
$ws.Cells.Item(27, 8).Value = "5.95"
$ws.Cells.Item(27, 9).Value = ""
$ws.Cells.Item(27, 10).Value = ""
$ws.Cells.Item(27, 11).Value = ""
$ws.Cells.Item(27, 12).Value = ""
$ws.Cells.Item(27, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
$ws.Cells.Item(27, 14).Value = "Tela Viva Haushaltspapier 3-lagig 4 Rollen 5.95 Schweizer Franken"
$ws.Cells.Item(27, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(28, 1).NumberFormat = "@"
$ws.Cells.Item(28, 1).Value = "6800946"
$ws.Cells.Item(28, 2).Value = "Hipp Natural zart duftend 3x48 Stück"
$ws.Cells.Item(28, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/hipp-natural-zart-duftend-3x48-stueck/p/6800946"
$ws.Cells.Item(28, 4).Value = "144ST"
$ws.Cells.Item(28, 5).Value = 2
$ws.Cells.Item(28, 7).Value = "Hipp"
$ws.Cells.Item(28, 8).NumberFormat = "@"
$ws.Cells.Item(28, 8).Value = "6.65"
$ws.Cells.Item(28, 9).Value = "0.05/1ST"
$ws.Cells.Item(28, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(28, 11).NumberFormat = "@"
$ws.Cells.Item(28, 11).Value = "0.05"
$ws.Cells.Item(28, 12).Value = "1ST"
$ws.Cells.Item(28, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Cells.Item(28, 14).Value = "Hipp Natural zart duftend 3x48 Stück 33% Aktion 6.65 Schweizer Franken statt 9.95 Schweizer Franken"
$ws.Cells.Item(28, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(29, 5).Value = 1
$ws.Cells.Item(29, 6).Value = 5
$ws.Cells.Item(29, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(30, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(31, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(32, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(33, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(34, 1).NumberFormat = "@"
$ws.Cells.Item(34, 1).Value = "7041905"
$ws.Cells.Item(34, 2).Value = "Tela Taschentücher Spa 10x10 Stück"
$ws.Cells.Item(34, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/tela-taschentuecher-spa-10x10-stueck/p/7041905"
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = "Tela"
$ws.Cells.Item(34, 8).NumberFormat = "@"
$ws.Cells.Item(34, 8).Value = "2.55"
$ws.Cells.Item(34, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Cells.Item(34, 14).Value = "Tela Taschentücher Spa 10x10 Stück 2.55 Schweizer Franken"
$ws.Cells.Item(34, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(35, 1).NumberFormat = "@"
$ws.Cells.Item(35, 1).Value = "7041906"
$ws.Cells.Item(35, 2).Value = "Tela Taschentücher Spa Box 90Stück"
$ws.Cells.Item(35, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/tela-taschentuecher-spa-box-90stueck/p/7041906"
$ws.Cells.Item(35, 4).Value = "90ST"
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = "Tela"
$ws.Cells.Item(35, 8).NumberFormat = "@"
$ws.Cells.Item(35, 8).Value = "3.10"
$ws.Cells.Item(35, 9).Value = "0.03/1ST"
$ws.Cells.Item(35, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(35, 11).NumberFormat = "@"
$ws.Cells.Item(35, 11).Value = "0.03"
$ws.Cells.Item(35, 12).Value = "1ST"
$ws.Cells.Item(35, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Cells.Item(35, 14).Value = "Tela Taschentücher Spa Box 90Stück 3.10 Schweizer Franken"
$ws.Cells.Item(35, 15).Value = "2022-09-16 21:00:05"
$ws.Cells.Item(36, 1).NumberFormat = "@"
$ws.Cells.Item(36, 1).Value = "7042418"
$ws.Cells.Item(36, 2).Value = "Tempo feucht Limited Edition"
$ws.Cells.Item(36, 3).Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/tempo-feucht-limited-edition/p/7042418"
$ws.Cells.Item(36, 4).Value = "86ST"
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = "Tempo"
$ws.Cells.Item(36, 8).NumberFormat = "@"
$ws.Cells.Item(36, 8).Value = "6.95"
$ws.Cells.Item(36, 9).Value = "0.08/1ST"
$ws.Cells.Item(36, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(36, 11).NumberFormat = "@"
$ws.Cells.Item(36, 11).Value = "0.08"
$ws.Cells.Item(36, 12).Value = "1ST"
$ws.Cells.Item(36, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Cells.Item(36, 14).Value = "Tempo feucht Limited Edition 6.95 Schweizer Franken"
$ws.Cells.Item(36, 15).Value = "2022-09-16 21:00:05"
